$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.003.09"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.860.31"
$ws.Range("E3").Value = "  -0.91%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.16"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5141"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3839"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08233"
$ws.Range("E9").Value = "  -9.71%  "

$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.46"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.196"
$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.57"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("D14").Value = "1.868.07"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.266"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001096"
$ws.Range("E17").Value = "  -1.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.64"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06651"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.007"
$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").Value = "28.029.52"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -3.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.269"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("D26").Value = "2.075.38"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.512"
$ws.Range("E27").Value = "  -2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.65"
$ws.Range("E28").Value = "  +0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.47"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.67"
$ws.Range("E30").Value = "  -1.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1066"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.030"
$ws.Range("E32").Value = "  -3.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.913"
$ws.Range("E33").Value = "  +5.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.593"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.395"
$ws.Range("E35").Value = "  -3.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02412"
$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06501"
$ws.Range("E37").Value = "  -1.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2180"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6543"
$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.198"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.998"
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.18"
$ws.Range("E43").Value = "  -3.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6176"
$ws.Range("E44").Value = "  +2.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.04"
$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.284"
$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.669"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.007"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.81"
$ws.Range("E50").Value = "  -0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.36"
$ws.Range("E51").Value = "  -1.88%  "
